# Updates the "Pop" column (B2:B452) on Sheet1 with refreshed regression
# output values (Year stays the same in column A).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPopValues = @(
    397832.5590765636,  # row 2, Year 1600
    398296.4294421462,  # row 3, Year 1601
    398766.9494714359,  # row 4, Year 1602
    399244.2136033769,  # row 5, Year 1603
    399728.3175924707,  # row 6, Year 1604
    400219.3585263552,  # row 7, Year 1605
    400717.4348435986,  # row 8, Year 1606
    401222.6463517063,  # row 9, Year 1607
    401735.0942453461,  # row 10, Year 1608
    402254.8811247918,  # row 11, Year 1609
    402782.1110145871,  # row 12, Year 1610
    403316.8893824328,  # row 13, Year 1611
    403859.3231582976,  # row 14, Year 1612
    404409.5207537565,  # row 15, Year 1613
    404967.5920815563,  # row 16, Year 1614
    405533.6485754117,  # row 17, Year 1615
    406107.8032100335,  # row 18, Year 1616
    406690.1705213905,  # row 19, Year 1617
    407280.8666272069,  # row 20, Year 1618
    407880.0092476972,  # row 21, Year 1619
    408487.7177265403,  # row 22, Year 1620
    409104.113052095,  # row 23, Year 1621
    409729.3178788575,  # row 24, Year 1622
    410363.4565491638,  # row 25, Year 1623
    411006.6551151386,  # row 26, Year 1624
    411659.0413608913,  # row 27, Year 1625
    412320.7448249621,  # row 28, Year 1626
    412991.8968230191,  # row 29, Year 1627
    413672.6304708078,  # row 30, Year 1628
    414363.0807073556,  # row 31, Year 1629
    415063.3843184313,  # row 32, Year 1630
    415773.6799602626,  # row 33, Year 1631
    416494.1081835119,  # row 34, Year 1632
    417224.8114575126,  # row 35, Year 1633
    417965.9341947667,  # row 36, Year 1634
    418717.6227757058,  # row 37, Year 1635
    419480.0255737154,  # row 38, Year 1636
    420253.2929804255,  # row 39, Year 1637
    421037.5774312675,  # row 40, Year 1638
    421833.0334312984,  # row 41, Year 1639
    422639.817581295,  # row 42, Year 1640
    423458.0886041162,  # row 43, Year 1641
    424288.0073713376,  # row 44, Year 1642
    425129.7369301571,  # row 45, Year 1643
    425983.4425305732,  # row 46, Year 1644
    426849.2916528364,  # row 47, Year 1645
    427727.4540351748,  # row 48, Year 1646
    428618.1017017938,  # row 49, Year 1647
    429521.4089911515,  # row 50, Year 1648
    430437.5525845094,  # row 51, Year 1649
    431366.7115347585,  # row 52, Year 1650
    432309.0672955227,  # row 53, Year 1651
    433264.8037505372,  # row 54, Year 1652
    434234.1072433044,  # row 55, Year 1653
    435217.1666070256,  # row 56, Year 1654
    436214.1731948093,  # row 57, Year 1655
    437225.3209101555,  # row 58, Year 1656
    438250.8062377158,  # row 59, Year 1657
    439290.8282743281,  # row 60, Year 1658
    440345.5887603274,  # row 61, Year 1659
    441415.2921111293,  # row 62, Year 1660
    442500.1454490866,  # row 63, Year 1661
    443600.3586356198,  # row 64, Year 1662
    444716.1443036161,  # row 65, Year 1663
    445847.7178901008,  # row 66, Year 1664
    446995.2976691747,  # row 67, Year 1665
    448159.1047852209,  # row 68, Year 1666
    449339.3632863744,  # row 69, Year 1667
    450536.3001582564,  # row 70, Year 1668
    451750.1453579699,  # row 71, Year 1669
    452981.1318483533,  # row 72, Year 1670
    454229.4956324917,  # row 73, Year 1671
    455495.4757884815,  # row 74, Year 1672
    456779.3145044465,  # row 75, Year 1673
    458081.2571138024,  # row 76, Year 1674
    459401.552130766,  # row 77, Year 1675
    460740.4512861067,  # row 78, Year 1676
    462098.2095631349,  # row 79, Year 1677
    463475.0852339258,  # row 80, Year 1678
    464871.3398957726,  # row 81, Year 1679
    466287.2385078656,  # row 82, Year 1680
    467723.0494281924,  # row 83, Year 1681
    469179.0444506549,  # row 84, Year 1682
    470655.4988423969,  # row 85, Year 1683
    472152.6913813382,  # row 86, Year 1684
    473670.9043939086,  # row 87, Year 1685
    475210.4237929765,  # row 88, Year 1686
    476771.5391159654,  # row 89, Year 1687
    478354.5435631521,  # row 90, Year 1688
    479959.7340361397,  # row 91, Year 1689
    481587.411176498,  # row 92, Year 1690
    483237.8794045645,  # row 93, Year 1691
    484911.4469583969,  # row 94, Year 1692
    486608.4259328704,  # row 95, Year 1693
    488329.1323189103,  # row 96, Year 1694
    490073.8860428514,  # row 97, Year 1695
    491843.0110059141,  # row 98, Year 1696
    493636.8351237897,  # row 99, Year 1697
    495455.6903663218,  # row 100, Year 1698
    497299.912797276,  # row 101, Year 1699
    499169.8426141852,  # row 102, Year 1700
    501065.8241882601,  # row 103, Year 1701
    502988.2061043535,  # row 104, Year 1702
    504937.3412009649,  # row 105, Year 1703
    506913.5866102735,  # row 106, Year 1704
    508917.3037981869,  # row 107, Year 1705
    510948.8586043911,  # row 108, Year 1706
    513008.6212823879,  # row 109, Year 1707
    515096.9665395056,  # row 110, Year 1708
    517214.2735768665,  # row 111, Year 1709
    519360.9261292972,  # row 112, Year 1710
    521537.3125051646,  # row 113, Year 1711
    523743.8256261206,  # row 114, Year 1712
    525980.8630667382,  # row 115, Year 1713
    528248.8270940225,  # row 116, Year 1714
    530548.1247067748,  # row 117, Year 1715
    532879.1676747951,  # row 118, Year 1716
    535242.372577899,  # row 119, Year 1717
    537638.1608447314,  # row 120, Year 1718
    540066.9587913544,  # row 121, Year 1719
    542529.1976595897,  # row 122, Year 1720
    545025.3136550891,  # row 123, Year 1721
    547555.747985116,  # row 124, Year 1722
    550120.9468960072,  # row 125, Year 1723
    552721.3617102974,  # row 126, Year 1724
    555357.4488634759,  # row 127, Year 1725
    558029.6699403513,  # row 128, Year 1726
    560738.4917109993,  # row 129, Year 1727
    563484.3861662628,  # row 130, Year 1728
    566267.8305527789,  # row 131, Year 1729
    569089.3074075015,  # row 132, Year 1730
    571949.3045916918,  # row 133, Year 1731
    574848.3153243433,  # row 134, Year 1732
    577786.8382150128,  # row 135, Year 1733
    580765.3772960225,  # row 136, Year 1734
    583784.4420540005,  # row 137, Year 1735
    586844.5474607274,  # row 138, Year 1736
    589946.21400325,  # row 139, Year 1737
    593089.9677132313,  # row 140, Year 1738
    596276.3401954947,  # row 141, Year 1739
    599505.8686557283,  # row 142, Year 1740
    602779.0959273088,  # row 143, Year 1741
    606096.5704972056,  # row 144, Year 1742
    609458.846530925,  # row 145, Year 1743
    612866.4838964515,  # row 146, Year 1744
    616320.0481871443,  # row 147, Year 1745
    619820.1107435459,  # row 148, Year 1746
    623367.2486740567,  # row 149, Year 1747
    626962.0448744309,  # row 150, Year 1748
    630605.0880460463,  # row 151, Year 1749
    634296.9727129011,  # row 152, Year 1750
    638038.2992372868,  # row 153, Year 1751
    641829.6738340901,  # row 154, Year 1752
    645671.7085836699,  # row 155, Year 1753
    649565.0214432592,  # row 156, Year 1754
    653510.2362568378,  # row 157, Year 1755
    657507.9827634224,  # row 158, Year 1756
    661558.8966037163,  # row 159, Year 1757
    665663.6193250662,  # row 160, Year 1758
    669822.7983846645,  # row 161, Year 1759
    674037.08715094,  # row 162, Year 1760
    678307.1449030782,  # row 163, Year 1761
    682633.6368286088,  # row 164, Year 1762
    687017.2340189982,  # row 165, Year 1763
    691458.6134631842,  # row 166, Year 1764
    695958.4580389897,  # row 167, Year 1765
    700517.4565023483,  # row 168, Year 1766
    705136.3034742755,  # row 169, Year 1767
    709815.6994255192,  # row 170, Year 1768
    714556.3506588193,  # row 171, Year 1769
    719358.9692887075,  # row 172, Year 1770
    724224.2732187754,  # row 173, Year 1771
    729152.9861163417,  # row 174, Year 1772
    734145.8373844405,  # row 175, Year 1773
    739203.5621310624,  # row 176, Year 1774
    744326.9011355699,  # row 177, Year 1775
    749516.6008122107,  # row 178, Year 1776
    754773.4131706553,  # row 179, Year 1777
    760098.0957734771,  # row 180, Year 1778
    765491.4116904987,  # row 181, Year 1779
    770954.1294499234,  # row 182, Year 1780
    776487.0229861713,  # row 183, Year 1781
    782090.8715843403,  # row 184, Year 1782
    787766.4598212066,  # row 185, Year 1783
    793514.5775026841,  # row 186, Year 1784
    799336.0195976582,  # row 187, Year 1785
    805231.5861681097,  # row 188, Year 1786
    811202.0822954428,  # row 189, Year 1787
    817248.3180029325,  # row 190, Year 1788
    823371.1081742048,  # row 191, Year 1789
    829571.2724676628,  # row 192, Year 1790
    835849.6352267719,  # row 193, Year 1791
    842207.0253861155,  # row 194, Year 1792
    848644.2763731342,  # row 195, Year 1793
    855162.2260054604,  # row 196, Year 1794
    861761.7163837585,  # row 197, Year 1795
    868443.5937799852,  # row 198, Year 1796
    875208.7085209771,  # row 199, Year 1797
    882057.9148672806,  # row 200, Year 1798
    888992.0708871339,  # row 201, Year 1799
    896012.0383255137,  # row 202, Year 1800
    903118.6824681584,  # row 203, Year 1801
    910312.8720004798,  # row 204, Year 1802
    917595.4788612761,  # row 205, Year 1803
    924967.3780911621,  # row 206, Year 1804
    932429.4476756258,  # row 207, Year 1805
    939982.568382632,  # row 208, Year 1806
    947627.6235946855,  # row 209, Year 1807
    955365.4991352689,  # row 210, Year 1808
    963197.0830895789,  # row 211, Year 1809
    971123.2656194749,  # row 212, Year 1810
    979144.9387725629,  # row 213, Year 1811
    987262.9962853347,  # row 214, Year 1812
    995478.3333802922,  # row 215, Year 1813
    1003791.846556972,  # row 216, Year 1814
    1012204.43337681,  # row 217, Year 1815
    1020716.992241764,  # row 218, Year 1816
    1029330.422166637,  # row 219, Year 1817
    1038045.622545025,  # row 220, Year 1818
    1046863.492908833,  # row 221, Year 1819
    1055784.932681302,  # row 222, Year 1820
    1064810.840923472,  # row 223, Year 1821
    1073942.116074056,  # row 224, Year 1822
    1083179.65568264,  # row 225, Year 1823
    1092524.356136189,  # row 226, Year 1824
    1101977.112378796,  # row 227, Year 1825
    1111538.817624652,  # row 228, Year 1826
    1121210.363064175,  # row 229, Year 1827
    1130992.637563294,  # row 230, Year 1828
    1140886.527355836,  # row 231, Year 1829
    1150892.915729014,  # row 232, Year 1830
    1161012.682701977,  # row 233, Year 1831
    1171246.704697425,  # row 234, Year 1832
    1181595.854206273,  # row 235, Year 1833
    1192060.999445359,  # row 236, Year 1834
    1202643.004008207,  # row 237, Year 1835
    1213342.726508844,  # row 238, Year 1836
    1224161.020218697,  # row 239, Year 1837
    1235098.732696577,  # row 240, Year 1838
    1246156.70541179,  # row 241, Year 1839
    1257335.773360402,  # row 242, Year 1840
    1268636.764674702,  # row 243, Year 1841
    1280060.500225909,  # row 244, Year 1842
    1291607.793220183,  # row 245, Year 1843
    1303279.448787997,  # row 246, Year 1844
    1315076.26356695,  # row 247, Year 1845
    1326999.025278084,  # row 248, Year 1846
    1339048.512295814,  # row 249, Year 1847
    1351225.49321155,  # row 250, Year 1848
    1363530.726391121,  # row 251, Year 1849
    1375964.959526114,  # row 252, Year 1850
    1388528.929179256,  # row 253, Year 1851
    1401223.360323952,  # row 254, Year 1852
    1414048.965878149,  # row 255, Year 1853
    1427006.446232648,  # row 256, Year 1854
    1440096.488774045,  # row 257, Year 1855
    1453319.767402455,  # row 258, Year 1856
    1466676.94204422,  # row 259, Year 1857
    1480168.658159765,  # row 260, Year 1858
    1493795.546246829,  # row 261, Year 1859
    1507558.221339267,  # row 262, Year 1860
    1521457.282501651,  # row 263, Year 1861
    1535493.312319896,  # row 264, Year 1862
    1549666.876388175,  # row 265, Year 1863
    1563978.522792351,  # row 266, Year 1864
    1578428.78159021,  # row 267, Year 1865
    1593018.16428877,  # row 268, Year 1866
    1607747.163318953,  # row 269, Year 1867
    1622616.251507914,  # row 270, Year 1868
    1637625.881549345,  # row 271, Year 1869
    1652776.485472083,  # row 272, Year 1870
    1668068.474107332,  # row 273, Year 1871
    1683502.23655487,  # row 274, Year 1872
    1699078.139648581,  # row 275, Year 1873
    1714796.527421686,  # row 276, Year 1874
    1730657.720572045,  # row 277, Year 1875
    1746662.015927921,  # row 278, Year 1876
    1762809.685914616,  # row 279, Year 1877
    1779100.978022366,  # row 280, Year 1878
    1795536.114275941,  # row 281, Year 1879
    1812115.290706364,  # row 282, Year 1880
    1828838.676825195,  # row 283, Year 1881
    1845706.415101833,  # row 284, Year 1882
    1862718.620444292,  # row 285, Year 1883
    1879875.379683924,  # row 286, Year 1884
    1897176.751064565,  # row 287, Year 1885
    1914622.763736594,  # row 288, Year 1886
    1932213.417256388,  # row 289, Year 1887
    1949948.681091703,  # row 290, Year 1888
    1967828.494133458,  # row 291, Year 1889
    1985852.764214458,  # row 292, Year 1890
    2004021.367635589,  # row 293, Year 1891
    2022334.148699983,  # row 294, Year 1892
    2040790.919255728,  # row 295, Year 1893
    2059391.458247627,  # row 296, Year 1894
    2078135.511278572,  # row 297, Year 1895
    2097022.790181082,  # row 298, Year 1896
    2116052.972599539,  # row 299, Year 1897
    2135225.701583693,  # row 300, Year 1898
    2154540.585193988,  # row 301, Year 1899
    2173997.196119258,  # row 302, Year 1900
    2193595.071307354,  # row 303, Year 1901
    2213333.711609273,  # row 304, Year 1902
    2233212.581437316,  # row 305, Year 1903
    2253231.108437855,  # row 306, Year 1904
    2273388.683179252,  # row 307, Year 1905
    2293684.658855467,  # row 308, Year 1906
    2314118.351005916,  # row 309, Year 1907
    2334689.037252106,  # row 310, Year 1908
    2355395.957051586,  # row 311, Year 1909
    2376238.311469723,  # row 312, Year 1910
    2397215.262969852,  # row 313, Year 1911
    2418325.935222272,  # row 314, Year 1912
    2439569.412932614,  # row 315, Year 1913
    2460944.741690062,  # row 316, Year 1914
    2482450.927835892,  # row 317, Year 1915
    2504086.938352818,  # row 318, Year 1916
    2525851.70077557,  # row 319, Year 1917
    2547744.103123159,  # row 320, Year 1918
    2569762.993853247,  # row 321, Year 1919
    2591907.18183902,  # row 322, Year 1920
    2614175.436368962,  # row 323, Year 1921
    2636566.487169902,  # row 324, Year 1922
    2659079.02445367,  # row 325, Year 1923
    2681711.698987711,  # row 326, Year 1924
    2704463.12218997,  # row 327, Year 1925
    2727331.866248315,  # row 328, Year 1926
    2750316.464264795,  # row 329, Year 1927
    2773415.410424958,  # row 330, Year 1928
    2796627.160192458,  # row 331, Year 1929
    2819950.130529141,  # row 332, Year 1930
    2843382.700140788,  # row 333, Year 1931
    2866923.209748657,  # row 334, Year 1932
    2890569.962386938,  # row 335, Year 1933
    2914321.223726217,  # row 336, Year 1934
    2938175.222423006,  # row 337, Year 1935
    2962130.150495375,  # row 338, Year 1936
    2986184.163724699,  # row 339, Year 1937
    3010335.382083478,  # row 340, Year 1938
    3034581.89018919,  # row 341, Year 1939
    3058921.737784087,  # row 342, Year 1940
    3083352.940240819,  # row 343, Year 1941
    3107873.479093733,  # row 344, Year 1942
    3132481.302595685,  # row 345, Year 1943
    3157174.326300141,  # row 346, Year 1944
    3181950.433668325,  # row 347, Year 1945
    3206807.476701162,  # row 348, Year 1946
    3231743.276595684,  # row 349, Year 1947
    3256755.624425595,  # row 350, Year 1948
    3281842.281845596,  # row 351, Year 1949
    3307000.981819103,  # row 352, Year 1950
    3332229.429368904,  # row 353, Year 1951
    3357525.30235032,  # row 354, Year 1952
    3382886.252246355,  # row 355, Year 1953
    3408309.904984341,  # row 356, Year 1954
    3433793.861773507,  # row 357, Year 1955
    3459335.699962917,  # row 358, Year 1956
    3484932.973919145,  # row 359, Year 1957
    3510583.215923085,  # row 360, Year 1958
    3536283.937085204,  # row 361, Year 1959
    3562032.628278572,  # row 362, Year 1960
    3587826.76108896,  # row 363, Year 1961
    3613663.788781251,  # row 364, Year 1962
    3639541.14728142,  # row 365, Year 1963
    3665456.256173303,  # row 366, Year 1964
    3691406.519709335,  # row 367, Year 1965
    3717389.327834451,  # row 368, Year 1966
    3743402.057222301,  # row 369, Year 1967
    3769442.072322911,  # row 370, Year 1968
    3795506.726420919,  # row 371, Year 1969
    3821593.3627035,  # row 372, Year 1970
    3847699.315337046,  # row 373, Year 1971
    3873821.910551715,  # row 374, Year 1972
    3899958.467732885,  # row 375, Year 1973
    3926106.300518593,  # row 376, Year 1974
    3952262.717901987,  # row 377, Year 1975
    3978425.02533785,  # row 378, Year 1976
    4004590.525852221,  # row 379, Year 1977
    4030756.521154135,  # row 380, Year 1978
    4056920.312748532,  # row 381, Year 1979
    4083079.203049333,  # row 382, Year 1980
    4109230.49649172,  # row 383, Year 1981
    4135371.500642652,  # row 384, Year 1982
    4161499.527308633,  # row 385, Year 1983
    4187611.89363977,  # row 386, Year 1984
    4213705.923229167,  # row 387, Year 1985
    4239778.947206673,  # row 388, Year 1986
    4265828.305326081,  # row 389, Year 1987
    4291851.347044797,  # row 390, Year 1988
    4317845.432595082,  # row 391, Year 1989
    4343807.934045947,  # row 392, Year 1990
    4369736.236354792,  # row 393, Year 1991
    4395627.738407924,  # row 394, Year 1992
    4421479.854049065,  # row 395, Year 1993
    4447290.013095021,  # row 396, Year 1994
    4473055.662337665,  # row 397, Year 1995
    4498774.266531428,  # row 398, Year 1996
    4524443.309365517,  # row 399, Year 1997
    4550060.294420073,  # row 400, Year 1998
    4575622.746105551,  # row 401, Year 1999
    4601128.210584565,  # row 402, Year 2000
    4626574.25667555,  # row 403, Year 2001
    4651958.476737518,  # row 404, Year 2002
    4677278.487535312,  # row 405, Year 2003
    4702531.931084711,  # row 406, Year 2004
    4727716.47547681,  # row 407, Year 2005
    4752829.815681121,  # row 408, Year 2006
    4777869.674326848,  # row 409, Year 2007
    4802833.80246186,  # row 410, Year 2008
    4827719.980288858,  # row 411, Year 2009
    4852526.017878346,  # row 412, Year 2010
    4877249.755857934,  # row 413, Year 2011
    4901889.066077656,  # row 414, Year 2012
    4926441.852250914,  # row 415, Year 2013
    4950906.050570753,  # row 416, Year 2014
    4975279.630301181,  # row 417, Year 2015
    4999560.594343274,  # row 418, Year 2016
    5023746.979775855,  # row 419, Year 2017
    5047836.858370554,  # row 420, Year 2018
    5071828.337081086,  # row 421, Year 2019
    5095719.558506633,  # row 422, Year 2020
    5119508.701329217,  # row 423, Year 2021
    5143193.980725012,  # row 424, Year 2022
    5166773.64874955,  # row 425, Year 2023
    5190245.994696828,  # row 426, Year 2024
    5213609.345432317,  # row 427, Year 2025
    5236862.065699951,  # row 428, Year 2026
    5260002.558403155,  # row 429, Year 2027
    5283029.264860035,  # row 430, Year 2028
    5305940.665032855,  # row 431, Year 2029
    5328735.277731963,  # row 432, Year 2030
    5351411.660794363,  # row 433, Year 2031
    5373968.411237128,  # row 434, Year 2032
    5396404.165385907,  # row 435, Year 2033
    5418717.59897877,  # row 436, Year 2034
    5440907.427245693,  # row 437, Year 2035
    5462972.404963958,  # row 438, Year 2036
    5484911.326489818,  # row 439, Year 2037
    5506723.025766764,  # row 440, Year 2038
    5528406.37631075,  # row 441, Year 2039
    5549960.291172772,  # row 442, Year 2040
    5571383.722879186,  # row 443, Year 2041
    5592675.663350195,  # row 444, Year 2042
    5613835.143796923,  # row 445, Year 2043
    5634861.234597531,  # row 446, Year 2044
    5655753.045152843,  # row 447, Year 2045
    5676509.723721922,  # row 448, Year 2046
    5697130.457238129,  # row 449, Year 2047
    5717614.471106113,  # row 450, Year 2048
    5737961.02898027,  # row 451, Year 2049
    5758169.432525172  # row 452, Year 2050
)

$firstRow = 2
for ($i = 0; $i -lt $newPopValues.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 2).Value = $newPopValues[$i]
}
